$d = $word.ActiveDocument

# Locate the paragraph that ends the "References still work..." section
# (the BodyText paragraph immediately before the "Figures inside other
# environments" Heading1), so the four new example/colouroff paragraphs
# can be inserted right after it, and before the upcoming heading.
$anchorText = "References still work in the same way as in Bookdown. Now go to theorem 2.1 or proposition 2.2."

$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the anchor paragraph for the insertion point."
}

$anchorRange = $d.Content
$anchorRange.Start = $find.Parent.Start
$anchorRange.End = $find.Parent.End

# Build the four new BodyText paragraphs. Using a single InsertAfter call
# with embedded carriage returns makes each chunk its own new paragraph,
# and every new paragraph inherits the BodyText style from the paragraph
# it follows (matching the target formatting).
$newText = "`r" + `
    "Here is some text which is not part of the below example. " + "`r" + `
    " Examples: " + "`r" + `
    "You can turn off the colour and padding in html, ePub and Word for any newtheorem or inbuilt theorem type. You do this in the _bookdown.yml file by adding the theorem name to the colouroff style_with list." + "`r" + `
    "Here is some text which is not part of the above example."

$anchorRange.InsertAfter($newText)

Write-Output "Inserted colouroff example paragraphs."
